$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '35.670.24'
Set-TextCell 2 5 '  -2.84%  '

# Row 3
Set-TextCell 3 4 '1.982.25'
Set-TextCell 3 5 '  -3.81%  '

# Row 4
Set-TextCell 4 5 '  +0.08%  '

# Row 5
Set-TextCell 5 4 '245.71'
Set-TextCell 5 5 '  +0.57%  '

# Row 6
Set-TextCell 6 5 '  -4.60%  '

# Row 7
Set-TextCell 7 4 '59.52'
Set-TextCell 7 5 '  +7.40%  '

# Row 8
Set-TextCell 8 5 '  +0.07%  '

# Row 9
Set-TextCell 9 4 '59.02'
Set-TextCell 9 5 '  -1.21%  '

# Row 10
Set-TextCell 10 4 '0.364'
Set-TextCell 10 5 '  -0.80%  '

# Row 11
Set-TextCell 11 4 '0.0740'
Set-TextCell 11 5 '  -1.93%  '

# Row 12
Set-TextCell 12 5 '  -2.62%  '

# Row 13
Set-TextCell 13 4 '0.949'
Set-TextCell 13 5 '  +0.86%  '

# Row 14
Set-TextCell 14 5 '  -1.42%  '

# Row 15
Set-TextCell 15 4 '2.271.81'
Set-TextCell 15 5 '  -3.77%  '

# Row 16
Set-TextCell 16 4 '5.32'
Set-TextCell 16 5 '  -3.05%  '

# Row 17
Set-TextCell 17 2 'WrappedEther'
Set-TextCell 17 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 17 4 '2.010.45'
Set-TextCell 17 5 '  -2.43%  '

# Row 18
Set-TextCell 18 2 'Avalanche'
Set-TextCell 18 3 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 18 4 '18.75'
Set-TextCell 18 5 '  +9.01%  '

# Row 19
Set-TextCell 19 4 '35.597.29'
Set-TextCell 19 5 '  -2.83%  '

# Row 20
Set-TextCell 20 4 '71.65'
Set-TextCell 20 5 '  -0.82%  '

# Row 21
Set-TextCell 21 5 '  -2.03%  '

# Row 22
Set-TextCell 22 4 '5.22'
Set-TextCell 22 5 '  -1.38%  '

# Row 23
Set-TextCell 23 4 '233.35'

# Row 24
Set-TextCell 24 5 '  +0.21%  '

# Row 25
Set-TextCell 25 5 '  +17.03%  '

# Row 26
Set-TextCell 26 4 '2.27'
Set-TextCell 26 5 '  -4.63%  '

# Row 27
Set-TextCell 27 4 '165.25'
Set-TextCell 27 5 '  -0.08%  '

# Row 28
Set-TextCell 28 4 '9.21'
Set-TextCell 28 5 '  -1.63%  '

# Row 29
Set-TextCell 29 4 '19.28'
Set-TextCell 29 5 '  -4.76%  '

# Row 30
Set-TextCell 30 5 '  -2.65%  '

# Row 31
Set-TextCell 31 5 '  -4.65%  '

# Row 32
Set-TextCell 32 5 '  -6.34%  '

# Row 33
Set-TextCell 33 4 '0.0962'
Set-TextCell 33 5 '  +13.33%  '

# Row 34
Set-TextCell 34 4 '0.0598'
Set-TextCell 34 5 '  -0.50%  '

# Row 35
Set-TextCell 35 4 '2.44'
Set-TextCell 35 5 '  +9.70%  '

# Row 36
Set-TextCell 36 5 '  -3.69%  '

# Row 37
Set-TextCell 37 5 '  +0.04%  '

# Row 38
Set-TextCell 38 4 '1.79'
Set-TextCell 38 5 '  -2.75%  '

# Row 39
Set-TextCell 39 4 '5.50'
Set-TextCell 39 5 '  +9.70%  '

# Row 40
Set-TextCell 40 4 '1.23'
Set-TextCell 40 5 '  -1.91%  '

# Row 41
Set-TextCell 41 5 '  -1.39%  '

# Row 42
Set-TextCell 42 5 '  -1.17%  '

# Row 43
Set-TextCell 43 2 'Cronos'
Set-TextCell 43 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 43 4 '0.0918'
Set-TextCell 43 5 '  +1.36%  '

# Row 44
Set-TextCell 44 2 'FraxShare'
Set-TextCell 44 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 44 4 '7.82'
Set-TextCell 44 5 '  +0.94%  '

# Row 45
Set-TextCell 45 2 'InjectiveProtocol'
Set-TextCell 45 3 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 45 4 '16.47'
Set-TextCell 45 5 '  +1.94%  '

# Row 46
Set-TextCell 46 5 '  -1.83%  '

# Row 47
Set-TextCell 47 4 '93.54'
Set-TextCell 47 5 '  -1.58%  '

# Row 48
Set-TextCell 48 4 '1.365.73'
Set-TextCell 48 5 '  -3.21%  '

# Row 49
Set-TextCell 49 4 '2.89'
Set-TextCell 49 5 '  -0.80%  '

# Row 50
Set-TextCell 50 4 '47.01'
Set-TextCell 50 5 '  +2.95%  '

# Row 51
Set-TextCell 51 4 '2.30'
Set-TextCell 51 5 '  -0.01%  '
